# Apply the edit described by the commit diff:
#  1. Insert a new slide ("冒泡排序的优化") at position 8, using the same
#     "Title and Content" layout used by the other content slides. This
#     pushes the existing last slide ("模板题") down to position 9.
#  2. Add a new (empty) text box to slide 1, the title slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. New slide: "冒泡排序的优化" (Bubble sort optimisation), inserted right
#    before the last slide (the "模板题" template-problem slide).
# ---------------------------------------------------------------------
$contentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(8, $contentLayout)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "冒泡排序的优化"
$newTitle.TextFrame.TextRange.LanguageID = "zh-CN"

$newBody = $newSlide.Shapes.Item(2)
$newBody.TextFrame.TextRange.Text = "在没有进行冒泡操作的时候，我们可以提前截止，节省时间"
$newBody.TextFrame.TextRange.LanguageID = "zh-CN"

# ---------------------------------------------------------------------
# 2. Slide 1: add an empty text box (same size/position as in the author's
#    edit) below the subtitle block.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$textBox = $slide1.Shapes.AddTextbox(1, 578.4, 386.35, 320.0, 29.0)
$textBox.Name = "文本框 3"
$textBox.Fill.Visible = 0
$textBox.TextFrame.WordWrap = -1
$textBox.TextFrame.AutoSize = 1
# Re-assert the exact geometry, since AutoSize can nudge the height.
$textBox.Left = 578.4
$textBox.Top = 386.35
$textBox.Width = 320.0
$textBox.Height = 29.0

Write-Host "Slides: " $p.Slides.Count
